# Entregables/ChangeLog.xlsx - add "Primer informe" sheet after "Informe Inicial"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the selection on the original sheet before it loses focus
$ws1.Range("E9").Select() | Out-Null

# Duplicate the first sheet (keeps identical column widths / styles / merged cells)
# and place the copy right after it.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Primer informe"

# Update the title row
$ws2.Range("B2").Value = "PRIMER INFORME"

# Update the log entries
$ws2.Range("B5").Value = 43204
$ws2.Range("B6").Value = 43204
$ws2.Range("C6").Value = "elaboracion de la estructura e ideas principales"

# Clear out the remaining template rows (keep date-column formatting, drop the rest)
$ws2.Range("B7:B13").ClearContents()
$ws2.Range("C7:C13").Clear()

# Match the saved selection/active cell on the new sheet
$ws2.Range("C7").Select() | Out-Null
